{"js": "// The document body contains a single one-column table where each row\n// holds one benchmark metric value (the last three rows originally held\n// several tab-separated values packed into one cell/run).\n// We address cells by (row, column) index and replace each cell's whole\n// paragraph range text in place so the existing run formatting\n// (Times New Roman, sz 22) is preserved.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// row index (0-based) -> new cell text\nconst newValues = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"701\",\n  4: \"0.00001\",\n  5: \"0.00045\",\n  6: \"0.00011\",\n  8: \"0.00018\",\n  9: \"0.00021\",\n  10: \"0.00028\",\n  11: \"0.09287\",\n  43: \"99.94\",\n  44: \"0.09\",\n  45: \"143\",\n};\n\nconst rowIndices = Object.keys(newValues).map(Number);\n\n// Grab the first paragraph of each target cell.\nconst paragraphs = {};\nfor (const rowIndex of rowIndices) {\n  const cell = table.getCell(rowIndex, 0);\n  cell.body.paragraphs.load(\"items\");\n  paragraphs[rowIndex] = cell.body.paragraphs;\n}\nawait context.sync();\n\n// Replace the whole paragraph (and hence whole cell) text with the new\n// value. Using the paragraph Range keeps the original run's rPr\n// (font/size) instead of inserting a brand-new, unformatted run.\nfor (const rowIndex of rowIndices) {\n  const para = paragraphs[rowIndex].items[0];\n  const range = para.getRange();\n  range.insertText(newValues[rowIndex], Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The document body contains a single one-column table where each row\n# holds one benchmark metric value (the last three rows originally held\n# several tab-separated values packed into one cell).\n# Addressing cells by (row, column) and assigning Cell.Range.Text keeps\n# the existing run formatting (Times New Roman, sz 22) intact while\n# replacing the cell's whole content with the new value.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# 1-based row index -> new cell text\n$newValues = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"701\"\n    5  = \"0.00001\"\n    6  = \"0.00045\"\n    7  = \"0.00011\"\n    9  = \"0.00018\"\n    10 = \"0.00021\"\n    11 = \"0.00028\"\n    12 = \"0.09287\"\n    44 = \"99.94\"\n    45 = \"0.09\"\n    46 = \"143\"\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $cell = $tbl.Cell($rowIndex, 1)\n    $cell.Range.Text = $newValues[$rowIndex]\n}\n"}
